$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: unchanged title/description/expected text, add Pass/Fail result ---
$ws.Range("E3").Value = "Pass"

# --- Row 4: fix the Description wording (was referring to "Request" button,
# should refer to "Search" button) and add Pass/Fail result ---
$ws.Range("C4").Value = 'verify that clicking on "Search" button shows success message'
$ws.Range("E4").Value = "Pass"

# --- Row 5: add Pass/Fail result ---
$ws.Range("E5").Value = "Pass"

# --- Row 6: add Pass/Fail result ---
$ws.Range("E6").Value = "Pass"

# --- Row 7: reword from "Submit" button to "Request" button, add result ---
$ws.Range("C7").Value = 'verify that stakeholder clicking on "Request" button sends stakeholder information to the registrar'
$ws.Range("D7").Value = 'Clicking on "Request" button should send stakeholder information to the registrar'
$ws.Range("E7").Value = "Pass"

# --- Row 8: replace with the new "after clicking Request" test case, add a
# failing result plus a remark column ---
$ws.Range("C8").Value = 'verify that after clicking "Request" button stakeholder information send to UGC, Register, Program Officer of corresponding department and student'
$ws.Range("D8").Value = 'After clicking "Request" button stakeholder information should send to UGC, Register, Program Officer of corresponding department and student'
$ws.Range("E8").Value = "fail"
$ws.Range("F8").Value = "In student, he/she does not get any message"
$ws.Range("F8").VerticalAlignment = -4160
$ws.Rows(8).RowHeight = 60

# --- Row 9 (new): the old "link / student profile" test case, renumbered as
# test case 7, moved down from row 8 ---
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = "verify that stakeholder clicking on link shows student profile with certificate and other information"
$ws.Range("D9").Value = "Clicking on link should show student profile with certificate and other information"
$ws.Rows(9).RowHeight = 45

# The stray formatted-but-empty row 11 from the old sheet is gone in the new
# layout.
$ws.Rows(11).Delete()

# --- View state: active cell is now C9, and the view is scrolled so row 4 is
# the top visible row ---
$ws.Range("C9").Select()
$excel.ActiveWindow.ScrollRow = 4
